# Updated capital structure database
# Refresh the computed capital-structure metrics for the two Venezuela
# "Banks (Regional)" rows (row 2: industry aggregate, row 3: Banco
# Occidental de Descuento) with the latest source figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2, 3) {
    # Margin figures that collapse to (negative) zero under the new data.
    $ws.Range("G$r").Value = -0
    $ws.Range("H$r").Value = -0
    $ws.Range("I$r").Value = -0
    $ws.Range("J$r").Value = -0

    # Trailing net income / net margin.
    $ws.Range("K$r").Value = -9.470000000000001
    $ws.Range("L$r").Value = 1.22987012987013

    # Cash, cash/market-cap, ROE and related spreads.
    $ws.Range("U$r").Value = 39.8
    $ws.Range("V$r").Value = 11.95195195195195
    $ws.Range("W$r").Value = -0.0439239332096475
    $ws.Range("X$r").Value = 0.1387371770143146
    $ws.Range("Y$r").Value = -0.1826611102239621
    $ws.Range("Z$r").Value = -0.3948717948717949

    # Cost of capital / ROIC spread.
    $ws.Range("AB$r").Value = 0.1387371770143146
    $ws.Range("AC$r").Value = -0.1387371770143146

    # Net debt and leverage ratios.
    $ws.Range("AG$r").Value = -39.8
    $ws.Range("AJ$r").Value = 1.09130792432136
    $ws.Range("AK$r").Value = -0.3754716981132075
}

# Row-specific sign flips on the cash-returned/payout ratios.
$ws.Range("O2").Value = -0
$ws.Range("R2").Value = -0

$ws.Range("O3").Value = 0
$ws.Range("R3").Value = 0
